# Actualización 11 de Mayo - Mañana
$wb = $excel.ActiveWorkbook

# ---- Sheet "1er Parcial" ----
$ws1 = $wb.Worksheets.Item("1er Parcial")

# Row 8: Promedio (I8) updated
$ws1.Range("I8").Value = 8.300000000000001

# Row 9: Aprobados/Reprobados/etc recalculated
$ws1.Range("E9").Value = 22
$ws1.Range("F9").Value = 3
$ws1.Range("G9").Value = 88
$ws1.Range("H9").Value = 12
$ws1.Range("I9").Value = 7.7
$ws1.Range("J9").Value = 3
$ws1.Range("K9").Value = 12

# Row 10: Promedio (I10) updated
$ws1.Range("I10").Value = 7.5

# Row 12: Promedio (I12) updated
$ws1.Range("I12").Value = 7.6

# ---- Sheet "2o Parcial" ----
$ws2 = $wb.Worksheets.Item("2o Parcial")

# Row 7
$ws2.Range("E7").Value = 25
$ws2.Range("F7").Value = 14
$ws2.Range("G7").Value = 64.09999999999999
$ws2.Range("H7").Value = 35.9
$ws2.Range("I7").Value = 7.2
$ws2.Range("J7").Value = 14
$ws2.Range("K7").Value = 35.9

# Row 8
$ws2.Range("E8").Value = 34
$ws2.Range("F8").Value = 7
$ws2.Range("G8").Value = 82.93000000000001
$ws2.Range("H8").Value = 17.07
$ws2.Range("I8").Value = 8.199999999999999
$ws2.Range("J8").Value = 7
$ws2.Range("K8").Value = 17.07

# Row 9
$ws2.Range("E9").Value = 21
$ws2.Range("F9").Value = 4
$ws2.Range("G9").Value = 84
$ws2.Range("H9").Value = 16
$ws2.Range("I9").Value = 7.6
$ws2.Range("J9").Value = 4
$ws2.Range("K9").Value = 16

# Row 10
$ws2.Range("E10").Value = 30
$ws2.Range("F10").Value = 9
$ws2.Range("G10").Value = 76.92
$ws2.Range("H10").Value = 23.08
$ws2.Range("I10").Value = 7.6
$ws2.Range("J10").Value = 9
$ws2.Range("K10").Value = 23.08

# Row 11
$ws2.Range("E11").Value = 21
$ws2.Range("F11").Value = 14
$ws2.Range("G11").Value = 60
$ws2.Range("H11").Value = 40
$ws2.Range("I11").Value = 7.6
$ws2.Range("J11").Value = 14
$ws2.Range("K11").Value = 40

# Row 12
$ws2.Range("E12").Value = 25
$ws2.Range("F12").Value = 11
$ws2.Range("G12").Value = 69.44
$ws2.Range("H12").Value = 30.56
$ws2.Range("I12").Value = 7.6
$ws2.Range("J12").Value = 11
$ws2.Range("K12").Value = 30.56
